# Add zero counts to test cases for prediction
# Appends 10 new "hand count" rows (camera_id / handcount=0 / diet=h2 / scorer=EK)
# to the bottom of Sheet1, mirroring the new test-case images that still
# need to be scored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$images = @(
    "IMG_5548.JPG",
    "IMG_5273.JPG",
    "IMG_5323.JPG",
    "IMG_5432.JPG",
    "IMG_5154.JPG",
    "IMG_5460.JPG",
    "IMG_4357.JPG",
    "IMG_5151.JPG",
    "IMG_5587.JPG",
    "IMG_5488.JPG"
)

$startRow = 185

# Pass 1: write the camera_id filenames (column A) for every new row first,
# so the new shared-string entries land in row order before anything else
# introduces a new unique string.
for ($i = 0; $i -lt $images.Length; $i++) {
    $row = $startRow + $i

    # Column A: camera_id (filename of the new test image). Written with no
    # inherited column style, matching the plain literal used for these rows.
    $ws.Cells.Item($row, 1).Value = $images[$i]
    $ws.Cells.Item($row, 1).Style = "Normal"
}

# Pass 2: fill in the rest of each row (handcount / diet / scorer).
for ($i = 0; $i -lt $images.Length; $i++) {
    $row = $startRow + $i

    # Column G: handcount placeholder of 0 for the not-yet-scored test case.
    $ws.Cells.Item($row, 7).Value = 0

    # Column H: diet/category code shared with the rest of the sheet.
    $ws.Cells.Item($row, 8).Value = "h2"

    # Column I: scorer initials for this batch (new shared string, added
    # once all the filenames above already claimed their slots).
    $ws.Cells.Item($row, 9).Value = "EK"
}

# Put the active selection where the author left off after entering the data.
[void]$ws.Range("F193").Select()
